# Append a new "Bag" method run (row 3) to the CELG stock-predictor data
# sheet, mirroring the extra sample captured once the VPN/stream-reader
# timeout & exception handling were added upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 2's formatting down to row 3 first (keeps the date style (s="1")
# on column A without introducing a brand new number-format/style entry).
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

# New data row.
$ws.Range("A3").Value = 42600.881053240744
$ws.Range("B3").Value = "Bag"
$ws.Range("C3").Value = 12323
$ws.Range("D3").Value = 12371
$ws.Range("E3").Value = 1465
$ws.Range("F3").Value = 169
$ws.Range("G3").Value = 112
$ws.Range("H3").Value = 59
$ws.Range("I3").Value = 39
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = 0

# Column A widened slightly (best-fit recalculated) to fit the new values.
$ws.Columns.Item(1).ColumnWidth = 14

Write-Host "Added row 3 (Bag sample) to WorkSheet 1"
